$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Designator text fixes -------------------------------------------------
$ws.Range("B3").Value = "C1"
$ws.Range("B4").Value = "C2"
$ws.Range("B6").Value = "LinReg1,LinReg2,LinReg3"

# --- 2) Drop the "Add to BOM" column (G) content, keep per-cell formatting ---
$ws.Range("G1:G12").ClearContents()

# --- 3) Re-colour the JLCPCB Part # column (F) from explicit grey to the ----
#        automatic/theme text colour, same font/size/alignment otherwise -----
$fRange = $ws.Range("F2:F12")
$fRange.Font.Name = "Helvetica Neue"
$fRange.Font.Size = 10
$fRange.Font.ThemeColor = 1
$fRange.HorizontalAlignment = -4108
$fRange.VerticalAlignment = -4108

# --- 4) Append two new BOM rows ----------------------------------------------
$ws.Range("A13").Value = 1
$ws.Range("B13").Value = "C12,C10,C14,C8,C4"
$ws.Range("C13").Value = "C_0603_1608Metric"
$ws.Range("D13").Value = 5
$ws.Range("E13").Value = "100uF"
$ws.Range("F13").Value = "C19702"

$ws.Range("A14").Value = 7
$ws.Range("B14").Value = "BConv1"
$ws.Range("C14").Value = "TPSM83100SIUR"
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = "TPSM83100SIUR"
$ws.Range("F14").Value = "C20346010"

# Formatting for A13:E14 matches rows 2-5 (A/B/C = style 1, D/E = style 2)
# -- copy it down from row 2.
$ws.Range("A2:E2").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A2:E2").Copy()
$ws.Range("A14").PasteSpecial(-4122)

# F13/F14 share the new automatic-colour font, left at the default (general)
# alignment -- i.e. no explicit horizontal/vertical alignment is applied.
$fNewRange = $ws.Range("F13:F14")
$fNewRange.Font.Name = "Helvetica Neue"
$fNewRange.Font.Size = 10
$fNewRange.Font.ThemeColor = 1

# Stray formatted-but-empty cells left next to the new row 14
$ws.Range("H14").Font.Name = "Helvetica Neue"
$ws.Range("H14").Font.Size = 10
$ws.Range("H14").Font.Color = 0
$ws.Range("H14").HorizontalAlignment = -4108
$ws.Range("H14").VerticalAlignment = -4108

$ws.Range("I14").HorizontalAlignment = -4131
$ws.Range("I14").VerticalAlignment = -4108

# --- 5) Trailing blank / spacer rows -----------------------------------------
$ws.Range("B22:D22").Value = ""
$ws.Range("B22:D22").HorizontalAlignment = -4108
$ws.Range("B22:D22").VerticalAlignment = -4108

$ws.Range("E22:F22").Font.Name = "Helvetica Neue"
$ws.Range("E22:F22").Font.Size = 10
$ws.Range("E22:F22").Font.Color = 0
$ws.Range("E22:F22").HorizontalAlignment = -4108
$ws.Range("E22:F22").VerticalAlignment = -4108

$ws.Range("G22").Font.Name = "Calibri"
$ws.Range("G22").Font.Size = 10
$ws.Range("G22").Font.ThemeColor = 1
$ws.Range("G22").HorizontalAlignment = -4108
$ws.Range("G22").VerticalAlignment = -4108

$ws.Range("G23").Font.Name = "Helvetica Neue"
$ws.Range("G23").Font.Size = 14
$ws.Range("G23").Font.Color = 16752192
$ws.Rows.Item(23).RowHeight = 18

# --- 6) Column widths to fit the new, wider designator/footprint text -------
$ws.Columns.Item(2).ColumnWidth = 20.83
$ws.Columns.Item(3).ColumnWidth = 17.67
$ws.Columns.Item(4).ColumnWidth = 7.33
$ws.Columns.Item(5).ColumnWidth = 15.17

# --- 7) Selection parks on the new spacer area, matching the saved file -----
$ws.Range("E19").Select()
